$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Arveja Verde" (Femacal de La Calera),
# dated 2022-03-17 (serial 44637). It belongs chronologically right before the
# existing row 50, so insert a fresh row there and push the rest of the table
# (old rows 50-57) down to 51-58.
$ws.Rows(50).Insert()

$newRow = 50
$ws.Cells.Item($newRow, 1).Value = 3
$ws.Cells.Item($newRow, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 44637
$ws.Cells.Item($newRow, 5).Value = 5
$ws.Cells.Item($newRow, 6).Value = 100112022
$ws.Cells.Item($newRow, 7).Value = "Arveja Verde"
$ws.Cells.Item($newRow, 8).Value = "Perfection"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 38
$ws.Cells.Item($newRow, 11).Value = 25000
$ws.Cells.Item($newRow, 12).Value = 25000
$ws.Cells.Item($newRow, 13).Value = 25000
$ws.Cells.Item($newRow, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item($newRow, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($newRow, 16).Value = 1000
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
